$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Credentials")

# --- Extend Credentials sheet with additional password test data ---
$ws1.Range("C1").Value = "Wrong Passwords"

$ws1.Range("C2").Value = 1234
$ws1.Range("C2").HorizontalAlignment = -4131  # xlLeft
$ws1.Range("D2").Value = "Lowercase"

$ws1.Range("C3").Value = "a1234"
$ws1.Range("D3").Value = "Capital and Minimum 8 Characters"

$ws1.Range("C4").Value = "'@@@"
$ws1.Range("D4").Value = "Not Satisfy any rule"

$ws1.Range("D1").Value = "Not satisfy"

$ws1.Columns.Item(3).ColumnWidth = 21.66666666666667
$ws1.Columns.Item(4).ColumnWidth = 30
$ws1.Columns.Item(5).ColumnWidth = 24.16666666666667

$ws1.Range("C4").Select()

# --- Add a new "Properties" worksheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Properties"

$ws2.Range("A1").Value = "Expected Color"
$ws2.Range("A2").Value = "rgba(255, 0, 0, 1)"
$ws2.Columns.Item(1).ColumnWidth = 28.5
$ws2.Range("A2").Select()

$ws1.Select()
